$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 591.25
$ws.Range("I28").Value = 648.1429000000001
$ws.Range("J28").Value = 193
$ws.Range("K28").Value = 648.1429000000001
$ws.Range("L28").Value = 193
$ws.Range("M28").Value = -163.1429000000001
$ws.Range("N28").Value = -1163
# row 40
$ws.Range("H40").Value = 4179.4
$ws.Range("I40").Value = 2240.8333
$ws.Range("J40").Value = 5471.778
$ws.Range("K40").Value = 2240.8333
$ws.Range("L40").Value = 5471.778
$ws.Range("M40").Value = -2065.8333
$ws.Range("N40").Value = -5821.778
# row 70
$ws.Range("H70").Value = 2050
$ws.Range("I70").Value = 1900
$ws.Range("K70").Value = 5700
$ws.Range("M70").Value = -5430
# row 73
$ws.Range("H73").Value = 2050
$ws.Range("I73").Value = 1900
$ws.Range("K73").Value = 5700
$ws.Range("M73").Value = -4764
# row 80
$ws.Range("H80").Value = 4252.75
$ws.Range("I80").Value = 1966.6666
$ws.Range("K80").Value = 5899.9998
$ws.Range("M80").Value = -4901.9998
# row 83
$ws.Range("H83").Value = 4252.75
$ws.Range("I83").Value = 1966.6666
$ws.Range("K83").Value = 17699.9994
$ws.Range("M83").Value = -12707.9994
# row 98
$ws.Range("H98").Value = 2000
$ws.Range("I98").Value = 2000
$ws.Range("K98").Value = 2000
$ws.Range("M98").Value = -502
# row 113
$ws.Range("H113").Value = 2900
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 354
$ws.Range("N113").ClearContents()
# row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
# row 132
$ws.Range("H132").Value = 1392.9231
$ws.Range("I132").Value = 1175.6666
$ws.Range("K132").Value = 3526.9998
$ws.Range("M132").Value = -996.9998000000001
# row 137
$ws.Range("H137").Value = 4328.75
$ws.Range("I137").Value = 4308.6
$ws.Range("J137").Value = 4362.3335
$ws.Range("K137").Value = 12925.8
$ws.Range("L137").Value = 13087.0005
$ws.Range("M137").Value = -10375.8
$ws.Range("N137").Value = -18187.0005

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 897
$ws.Range("I45").Value = 897
$ws.Range("K45").Value = 897
$ws.Range("M45").Value = -520
# row 61
$ws.Range("H61").Value = 2221.1177
$ws.Range("I61").Value = 1355
$ws.Range("K61").Value = 1355
$ws.Range("M61").Value = -1143
# row 102
$ws.Range("H102").Value = 1999.5
$ws.Range("I102").Value = 1999.5
$ws.Range("K102").Value = 1999.5
$ws.Range("M102").Value = -377.5
# row 110
$ws.Range("H110").Value = 1121.75
$ws.Range("I110").Value = 1121.75
$ws.Range("K110").Value = 1121.75
$ws.Range("M110").Value = 923.25
# row 122
$ws.Range("H122").Value = 5925.6665
$ws.Range("I122").Value = 6234.222
$ws.Range("K122").Value = 18702.666
$ws.Range("M122").Value = -16252.666
# row 135
$ws.Range("H135").Value = 1366666.6
$ws.Range("J135").Value = 1366666.6
$ws.Range("L135").Value = 1366666.6
$ws.Range("N135").Value = -1376806.6
# row 136
$ws.Range("H136").Value = 2221.1177
$ws.Range("I136").Value = 1355
$ws.Range("K136").Value = 4065
$ws.Range("M136").Value = -1515
# row 139
$ws.Range("H139").Value = 49990
$ws.Range("J139").Value = 49990
$ws.Range("L139").Value = 49990
$ws.Range("N139").Value = -60270

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 12249.5
$ws.Range("I134").Value = 1499
$ws.Range("K134").Value = 4497
$ws.Range("M134").Value = -1962

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 8
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# row 16
$ws.Range("H16").Value = 3299.8
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 3749.5
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3749.5
$ws.Range("M16").Value = -2713
$ws.Range("N16").Value = -4323.5
# row 31
$ws.Range("H31").Value = 1920.3478
$ws.Range("I31").Value = 1032.0625
$ws.Range("J31").Value = 3950.7144
$ws.Range("K31").Value = 1032.0625
$ws.Range("L31").Value = 3950.7144
$ws.Range("M31").Value = -737.0625
$ws.Range("N31").Value = -4540.7144
# row 34
$ws.Range("H34").Value = 1920.3478
$ws.Range("I34").Value = 1032.0625
$ws.Range("J34").Value = 3950.7144
$ws.Range("K34").Value = 1032.0625
$ws.Range("L34").Value = 3950.7144
$ws.Range("M34").Value = -830.0625
$ws.Range("N34").Value = -4354.7144
# row 113
$ws.Range("H113").Value = 3299.8
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3749.5
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8089.5
# row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
# row 134
$ws.Range("H134").Value = 2998.1428
$ws.Range("I134").Value = 2998.1428
$ws.Range("K134").Value = 8994.428400000001
$ws.Range("M134").Value = -6459.428400000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 92
$ws.Range("H92").Value = 900
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196
$ws.Range("M92").ClearContents()
# row 98
$ws.Range("H98").Value = 722
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 3000
$ws.Range("N98").Value = -5996
# row 107
$ws.Range("H107").Value = 809.6667
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 791.6
$ws.Range("K107").Value = 2700
$ws.Range("L107").Value = 2374.8
$ws.Range("M107").Value = -780
$ws.Range("N107").Value = -6214.8

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 633.55554
$ws.Range("I97").Value = 734.4286
$ws.Range("K97").Value = 734.4286
$ws.Range("M97").Value = -238.4286
# row 113
$ws.Range("H113").Value = 1300.3334
$ws.Range("I113").Value = 1254.6666
$ws.Range("K113").Value = 1254.6666
$ws.Range("M113").Value = 915.3334
# row 122
$ws.Range("H122").Value = 3885.2856
$ws.Range("I122").Value = 3699.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11098.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8648.5
$ws.Range("N122").Value = -19900
# row 126
$ws.Range("H126").Value = 1995
$ws.Range("J126").Value = 1995
$ws.Range("L126").Value = 5985
$ws.Range("N126").Value = -10925
# row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 54
$ws.Range("H54").Value = 26950
$ws.Range("J54").Value = 26950
$ws.Range("L54").Value = 26950
$ws.Range("N54").Value = -27990
# row 76
$ws.Range("H76").Value = 49900
$ws.Range("J76").Value = 49900
$ws.Range("L76").Value = 49900
$ws.Range("N76").Value = -50530
# row 79
$ws.Range("H79").Value = 49900
$ws.Range("J79").Value = 49900
$ws.Range("L79").Value = 49900
$ws.Range("N79").Value = -52084
# row 126
$ws.Range("H126").Value = 1487.625
$ws.Range("I126").Value = 779.6
$ws.Range("K126").Value = 2338.8
$ws.Range("M126").Value = 131.1999999999998
# row 132
$ws.Range("H132").Value = 5098.8
$ws.Range("I132").Value = 3298.8
$ws.Range("K132").Value = 9896.400000000001
$ws.Range("M132").Value = -7366.400000000001
# row 136
$ws.Range("H136").Value = 1660.6154
$ws.Range("I136").Value = 1153.4546
$ws.Range("J136").Value = 4450
$ws.Range("K136").Value = 3460.3638
$ws.Range("L136").Value = 13350
$ws.Range("M136").Value = -910.3638000000001
$ws.Range("N136").Value = -18450
